# Modified test data to include files in different index
#
# - "Sequence file" sheet gains 3 new rows (12-14) duplicating the existing
#   "cell_ID_2" / lane-001 block (rows 9-11) but for lane 002 (LANE INDEX 2),
#   with new file names for the L002 lane.
# - The active sheet/tab moves from "Sequence file" to "Donor organism".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Sequence file": duplicate rows 9:11 (the lane-001 cell_ID_2 block)
#    into rows 12:14 for the new lane-002 files.
# ---------------------------------------------------------------------
$ws8 = $wb.Worksheets.Item("Sequence file")

$srcRows = $ws8.Range("A9:O11")
$dstRows = $ws8.Range("A12:O14")
$srcRows.Copy($dstRows)
$ws8.Range("A12:O14").RowHeight = 16

$ws8.Range("A12").Value = "small_WSSS_THYst9384954_S1_L002_R1_002.fastq.gz"
$ws8.Range("A13").Value = "small_WSSS_THYst9384954_S1_L002_R2_002.fastq.gz"
$ws8.Range("A14").Value = "small_WSSS_THYst9384954_S1_L002_I1_002.fastq.gz"

# LANE INDEX column -> these new rows are lane 2
$ws8.Range("H12:H14").Value = 2

# Leave the selection on the last-entered cell, like Excel would after typing.
[void]$ws8.Range("H14").Select()

# ---------------------------------------------------------------------
# 2. Move the active tab from "Sequence file" to "Donor organism".
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Donor organism")
$ws5.Activate()
